# Re-sort the weekly price rows (rows 3-7) by the "Fecha" (date) column,
# ascending, keeping all other column values attached to their row.
# (Row 2 already holds the earliest date and does not move.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    [PSCustomObject]@{ D = 45072; L = "Segunda"; M = 100; N = 16000; O = 16000; P = 16000; R = "Provincia de Chacabuco"; S = 889 }
    [PSCustomObject]@{ D = 45072; L = "Segunda"; M = 100; N = 17000; O = 17000; P = 17000; R = "Provincia de Limarí";    S = 944 }
    [PSCustomObject]@{ D = 44257; L = "Primera"; M = 100; N = 14000; O = 15000; P = 14500; R = "Región Metropolitana";  S = 806 }
    [PSCustomObject]@{ D = 44252; L = "Primera"; M = 120; N = 13000; O = 14000; P = 13500; R = "Región Metropolitana";  S = 750 }
    [PSCustomObject]@{ D = 44253; L = "Primera"; M = 160; N = 14000; O = 15000; P = 14500; R = "Región Metropolitana";  S = 806 }
)

$orderedRows = $rows | Sort-Object -Property D

$targetRow = 3
foreach ($data in $orderedRows) {
    $ws.Cells.Item($targetRow, 4).Value  = $data.D   # D: Fecha
    $ws.Cells.Item($targetRow, 12).Value = $data.L   # L: Calidad
    $ws.Cells.Item($targetRow, 13).Value = $data.M   # M: Volumen
    $ws.Cells.Item($targetRow, 14).Value = $data.N   # N: Precio mínimo
    $ws.Cells.Item($targetRow, 15).Value = $data.O   # O: Precio máximo
    $ws.Cells.Item($targetRow, 16).Value = $data.P   # P: Precio promedio ponderado
    $ws.Cells.Item($targetRow, 18).Value = $data.R   # R: Origen
    $ws.Cells.Item($targetRow, 19).Value = $data.S   # S: Precio $/Kg
    $targetRow = $targetRow + 1
}
